$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.725.07'
$ws.Range("E2").Value = '  +1.21%  '
$ws.Range("D3").Value = '1.645.64'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.532'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.61%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0891'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.93%  '
$ws.Range("D12").Value = '1.879.17'
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").Value = '1.644.75'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").Value = '27.696.62'
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = '0.0₃0726'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.35%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.70%  '
$ws.Range("E24").Value = '  -2.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("E26").Value = '  -1.75%  '
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("E28").Value = '  -1.03%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0486'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("E33").Value = '  +1.55%  '
$ws.Range("D34").Value = '1.445.88'
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  +2.09%  '
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.902'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.36%  '
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.74%  '
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.57%  '
$ws.Range("D47").Value = '1.788.06'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("D50").Value = '0.0₆0108'
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("E51").Value = '  -1.86%  '
